$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Time" column (G) with header + values for the two data rows.
$ws.Range("G1").Value = "Time"
$ws.Range("G2").Value = "09:46:04 2024-05-14"
$ws.Range("G3").Value = "09:46:04 2024-05-14"

# Quantity correction for row 2.
$ws.Range("B2").Value = 15

# Order Number corrections (must keep leading zeros -> format column as text first).
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D2").Value = "000001323"
$ws.Range("D3").Value = "000001323"
